$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Expand the "Table13" structured table from B4:C11 to B4:C19
#    (8 new rows for use cases UC-08 .. UC-15)
# ---------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
for ($i = 0; $i -lt 8; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# ---------------------------------------------------------------
# 2. Fill in the new use-case IDs (column B, rows 12-19)
# ---------------------------------------------------------------
$ws.Range("B12").Value = "UC-08"
$ws.Range("B13").Value = "UC-09"
$ws.Range("B14").Value = "UC-10"
$ws.Range("B15").Value = "UC-11"
$ws.Range("B16").Value = "UC-12"
$ws.Range("B17").Value = "UC-13"
$ws.Range("B18").Value = "UC-14"
$ws.Range("B19").Value = "UC-15"

# ---------------------------------------------------------------
# 3. Fill in the use-case names (column C, rows 5-15); rows 16-19
#    are left blank (diagrams to be added later)
# ---------------------------------------------------------------
$ws.Range("C5").Value  = "Create Goods Received Note"
$ws.Range("C6").Value  = "Place Order"
$ws.Range("C7").Value  = "Make Payment"
$ws.Range("C8").Value  = "Track Order Status"
$ws.Range("C9").Value  = "Create Goods Delivery Note"
$ws.Range("C10").Value = "Update Payment Status"
$ws.Range("C11").Value = "View Stock Report"
$ws.Range("C12").Value = "View Revenue Report"
$ws.Range("C13").Value = "Manage User Accounts"
$ws.Range("C14").Value = "Manage Orders"
$ws.Range("C15").Value = "Manage Product Catalogue"

# ---------------------------------------------------------------
# 4. Formatting
# ---------------------------------------------------------------

# 4a. New ID cells (B12:B19) reuse the same style as the existing
#     ID column cells (B5:B11) - center/middle, wrap text.
$ws.Range("B11").Copy()
$ws.Range("B12:B19").PasteSpecial(-4122)

# 4b. Usecases column for the already-existing rows (C5:C11):
#     vertical-center + wrap text, general horizontal alignment.
$ws.Range("C5").HorizontalAlignment = 1
$ws.Range("C5").VerticalAlignment = -4108
$ws.Range("C5").WrapText = $true
$ws.Range("C5").Copy()
$ws.Range("C6:C11").PasteSpecial(-4122)

# 4c. Usecases column for the newly added rows (C12:C19): same
#     alignment as above.
$ws.Range("C12").Locked = $true
$ws.Range("C12").HorizontalAlignment = 1
$ws.Range("C12").VerticalAlignment = -4108
$ws.Range("C12").WrapText = $true
$ws.Range("C12").Copy()
$ws.Range("C13:C19").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------
# 5. View: scroll to show the new rows and move the selection
# ---------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C16").Select()

# ---------------------------------------------------------------
# 6. Page setup - portrait orientation
# ---------------------------------------------------------------
$ws.PageSetup.Orientation = 1
